# Updated cryptos list on Thu Apr 25 23:28:41 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    # Force the value to be stored as text (avoids Excel auto-converting
    # numeric-looking strings into floating point numbers), then restore
    # the cell's original (default) style so no formatting is introduced.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Cells.Item(2, 4) "64.652.13"
$ws.Cells.Item(2, 5).Value = "  +0.80%  "

Set-TextCell $ws.Cells.Item(3, 4) "3.162.16"
$ws.Cells.Item(3, 5).Value = "  +1.12%  "

$ws.Cells.Item(4, 5).Value = "  -0.09%  "

Set-TextCell $ws.Cells.Item(5, 4) "614.34"
$ws.Cells.Item(5, 5).Value = "  +1.03%  "

Set-TextCell $ws.Cells.Item(6, 4) "145.37"
$ws.Cells.Item(6, 5).Value = "  -0.99%  "

$ws.Cells.Item(7, 5).Value = "  -0.04%  "

Set-TextCell $ws.Cells.Item(8, 4) "3.159.44"
$ws.Cells.Item(8, 5).Value = "  +1.19%  "

Set-TextCell $ws.Cells.Item(9, 4) "0.525"
$ws.Cells.Item(9, 5).Value = "  +0.06%  "

$ws.Cells.Item(10, 5).Value = "  +0.86%  "

Set-TextCell $ws.Cells.Item(11, 4) "5.44"
$ws.Cells.Item(11, 5).Value = "  -1.38%  "

$ws.Cells.Item(12, 5).Value = "  -0.45%  "

$ws.Cells.Item(13, 5).Value = "  +1.15%  "

$ws.Cells.Item(14, 5).Value = "  -1.89%  "

Set-TextCell $ws.Cells.Item(15, 4) "3.677.99"
$ws.Cells.Item(15, 5).Value = "  +0.97%  "

$ws.Cells.Item(16, 5).Value = "  +3.49%  "

Set-TextCell $ws.Cells.Item(17, 4) "64.622.65"
$ws.Cells.Item(17, 5).Value = "  +0.68%  "

Set-TextCell $ws.Cells.Item(18, 4) "3.157.14"
$ws.Cells.Item(18, 5).Value = "  +1.38%  "

$ws.Cells.Item(19, 5).Value = "  -0.77%  "

Set-TextCell $ws.Cells.Item(20, 4) "479.46"
$ws.Cells.Item(20, 5).Value = "  +0.38%  "

$ws.Cells.Item(21, 5).Value = "  +0.94%  "

Set-TextCell $ws.Cells.Item(22, 4) "0.721"
$ws.Cells.Item(22, 5).Value = "  +2.72%  "

$ws.Cells.Item(23, 5).Value = "  +3.79%  "

Set-TextCell $ws.Cells.Item(24, 4) "13.78"
$ws.Cells.Item(24, 5).Value = "  +0.91%  "

Set-TextCell $ws.Cells.Item(25, 4) "83.91"
$ws.Cells.Item(25, 5).Value = "  +1.15%  "

$ws.Cells.Item(26, 5).Value = "  +0.12%  "

$ws.Cells.Item(27, 2).Value = "PancakeSwap"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell $ws.Cells.Item(27, 4) "2.81"
$ws.Cells.Item(27, 5).Value = "  -3.50%  "

$ws.Cells.Item(28, 2).Value = "RenderToken"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell $ws.Cells.Item(28, 4) "8.66"
$ws.Cells.Item(28, 5).Value = "  +3.17%  "

Set-TextCell $ws.Cells.Item(29, 4) "7.11"
$ws.Cells.Item(29, 5).Value = "  +5.60%  "

$ws.Cells.Item(30, 5).Value = "  -0.91%  "

$ws.Cells.Item(31, 5).Value = "  -5.03%  "

$ws.Cells.Item(32, 5).Value = "  +0.02%  "

$ws.Cells.Item(33, 5).Value = "  -0.83%  "

Set-TextCell $ws.Cells.Item(34, 4) "26.49"
$ws.Cells.Item(34, 5).Value = "  +0.65%  "

Set-TextCell $ws.Cells.Item(35, 4) "1.13"
$ws.Cells.Item(35, 5).Value = "  +2.50%  "

Set-TextCell $ws.Cells.Item(36, 4) "0.0₃0785"
$ws.Cells.Item(36, 5).Value = "  +8.96%  "

$ws.Cells.Item(37, 5).Value = "  -0.54%  "

Set-TextCell $ws.Cells.Item(38, 4) "53.21"
$ws.Cells.Item(38, 5).Value = "  -2.25%  "

$ws.Cells.Item(39, 5).Value = "  +3.99%  "

Set-TextCell $ws.Cells.Item(40, 4) "460.72"
$ws.Cells.Item(40, 5).Value = "  +2.57%  "

$ws.Cells.Item(41, 5).Value = "  +0.84%  "

$ws.Cells.Item(42, 5).Value = "  -2.33%  "

$ws.Cells.Item(43, 5).Value = "  -0.64%  "

Set-TextCell $ws.Cells.Item(44, 4) "2.859.69"
$ws.Cells.Item(44, 5).Value = "  +0.21%  "

$ws.Cells.Item(45, 5).Value = "  +3.10%  "

$ws.Cells.Item(46, 5).Value = "  -0.40%  "

$ws.Cells.Item(47, 5).Value = "  +6.53%  "

Set-TextCell $ws.Cells.Item(48, 4) "26.56"
$ws.Cells.Item(48, 5).Value = "  +0.85%  "

$ws.Cells.Item(49, 5).Value = "  +0.06%  "

Set-TextCell $ws.Cells.Item(50, 4) "35.57"
$ws.Cells.Item(50, 5).Value = "  +9.42%  "

$ws.Cells.Item(51, 5).Value = "  -0.07%  "
